$d = $word.ActiveDocument

# 1. Rename the TOC bookmark on the "Summary" heading paragraph.
#    Bookmark.Name is read-only in the Word object model, so recreate the
#    bookmark at the same range under the new name.
$oldBookmarkName = "_Toc194782465"
$newBookmarkName = "_Toc194871182"
if ($d.Bookmarks.Exists($oldBookmarkName)) {
    $bm = $d.Bookmarks.Item($oldBookmarkName)
    $bmRange = $bm.Range
    $bm.Delete()
    $d.Bookmarks.Add($newBookmarkName, $bmRange) | Out-Null
}

# 2. Add 1.15-line spacing (w:line="276" w:lineRule="auto") to the
#    paragraphs that do not already carry explicit line spacing.
$paras = $d.Paragraphs
$targetIndexes = @(2, 9, 10, 11, 12, 13)
foreach ($i in $targetIndexes) {
    $p = $paras.Item($i)
    $p.LineSpacingRule = 5   # wdLineSpaceMultiple
    $p.LineSpacing = 13.8    # 276 twentieths-of-a-point -> 1.15 * 12pt
}

# 3. Strip the direct paragraph formatting from the final paragraph and
#    trim its trailing space. Re-insert the run as clean OOXML so no
#    <w:pPr> survives on the paragraph.
$paras = $d.Paragraphs
$last = $paras.Item($paras.Count)
$lastRange = $last.Range
$cleanXml = '<w:p w14:paraId="785EA5F7" w14:textId="394ED8D8" w:rsidR="00BA28E4" w:rsidRPr="00305AF1" w:rsidRDefault="00305AF1" w:rsidP="00305AF1" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>Another integrated feature should allow users to set a reminder by selecting a date. The reminder will take input on the plant type and species (e.g., succulent, tropical, houseplant) and notify the user of when to next water their plant.</w:t></w:r></w:p>'
$lastRange.InsertXML($cleanXml) | Out-Null
